# Repull data, push all data, mean calculation
# Update column F ("dSF") values for the rows whose source data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 1
    6  = 2
    7  = -3
    8  = 3
    9  = 7
    10 = 4
    11 = -1
    12 = 3
    13 = -1
    14 = -5
    15 = -2
    16 = 2
    17 = -2
    18 = -2
    19 = 4
    20 = -2
    21 = -1
    22 = 0
    24 = 4
    25 = 3
    26 = -2
    27 = 2
    28 = -2
    29 = -2
    31 = 1
    32 = 7
    33 = -3
    34 = -1
    36 = -4
    38 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
